# Apply the betexplorer scrape-refresh edit:
#  - A handful of adjacent match rows had their F:V (match/odds) payload
#    re-ordered (rows shifted by one position as the source re-scraped),
#    while the leading A:E (index/country/league/season/kickoff) columns
#    stayed put.
#  - Two brand-new fixtures were appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Simple pairwise swaps of the F:V block between two adjacent rows.
# ---------------------------------------------------------------------
$swapPairs = @(
    @(8, 9),
    @(18, 19),
    @(29, 30),
    @(31, 32),
    @(42, 43),
    @(44, 45),
    @(51, 52),
    @(89, 90),
    @(94, 95),
    @(106, 107),
    @(123, 124),
    @(127, 128)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# ---------------------------------------------------------------------
# 2) Three-way rotation of the F:V block across rows 97, 98, 99:
#    new97 <- old99, new98 <- old97, new99 <- old98
# ---------------------------------------------------------------------
$row97 = $ws.Range("F97:V97")
$row98 = $ws.Range("F98:V98")
$row99 = $ws.Range("F99:V99")

$old97 = $row97.Value()
$old98 = $row98.Value()
$old99 = $row99.Value()

$row97.Value = $old99
$row98.Value = $old97
$row99.Value = $old98

# ---------------------------------------------------------------------
# 3) Append two new fixtures as rows 134 and 135 (after existing row 133).
#    Copy formats from row 133 first so styles (bold/bordered index column,
#    datetime-formatted kickoff column) match the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("A133:V133").Copy()
$ws.Range("A134:V135").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 134
$ws.Range("A134").Value = 133
$ws.Range("B134").Value = "turkey"
$ws.Range("C134").Value = "super-lig"
$ws.Range("D134").Value = "2023-2024"
$ws.Range("E134").Value = 45262.47916666666
$ws.Range("F134").Value = "Alanyaspor"
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = "Konyaspor"
$ws.Range("I134").Value = 2
$ws.Range("J134").Value = 2.38
$ws.Range("K134").Value = "25/11/2023 18:13"
$ws.Range("L134").Value = 2.7
$ws.Range("M134").Value = "02/12/2023 11:25"
$ws.Range("N134").Value = 3.46
$ws.Range("O134").Value = "25/11/2023 18:13"
$ws.Range("P134").Value = 3.41
$ws.Range("Q134").Value = "02/12/2023 11:21"
$ws.Range("R134").Value = 3.03
$ws.Range("S134").Value = "25/11/2023 18:13"
$ws.Range("T134").Value = 2.75
$ws.Range("U134").Value = "02/12/2023 11:27"
$ws.Range("V134").Value = "https://www.betexplorer.com/football/turkey/super-lig/alanyaspor-konyaspor/vXtWdcmo/"

# Row 135
$ws.Range("A135").Value = 134
$ws.Range("B135").Value = "turkey"
$ws.Range("C135").Value = "super-lig"
$ws.Range("D135").Value = "2023-2024"
$ws.Range("E135").Value = 45262.58333333334
$ws.Range("F135").Value = "Rizespor"
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = "Basaksehir"
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = 2.63
$ws.Range("K135").Value = "27/11/2023 18:12"
$ws.Range("L135").Value = 2.56
$ws.Range("M135").Value = "02/12/2023 13:17"
$ws.Range("N135").Value = 3.3
$ws.Range("O135").Value = "27/11/2023 18:12"
$ws.Range("P135").Value = 3.33
$ws.Range("Q135").Value = "02/12/2023 13:55"
$ws.Range("R135").Value = 2.87
$ws.Range("S135").Value = "27/11/2023 18:12"
$ws.Range("T135").Value = 2.97
$ws.Range("U135").Value = "02/12/2023 13:55"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/turkey/super-lig/rizespor-basaksehir/newSGIvH/"

Write-Output "done"
